$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 82, shifting existing data down.
$ws.Rows.Item(82).Insert()

# Populate the new row 82 with the new data point.
$ws.Range("A82").Value2 = 4
$ws.Range("B82").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C82").Value2 = "Los Lagos"
$ws.Range("D82").Value2 = 44494
$ws.Range("E82").Value2 = 10
$ws.Range("F82").Value2 = 100112037
$ws.Range("G82").Value2 = "Cebollín"
$ws.Range("H82").Value2 = "Sin especificar"
$ws.Range("I82").Value2 = "Primera"
$ws.Range("J82").Value2 = 80
$ws.Range("K82").Value2 = 5500
$ws.Range("L82").Value2 = 5500
$ws.Range("M82").Value2 = 5500
$ws.Range("N82").Value2 = "$/paquete 36 unidades"
$ws.Range("O82").Value2 = "Región Metropolitana"
$ws.Range("P82").Value2 = 153
$ws.Range("Q82").Value2 = 36
$ws.Range("R82").Value2 = "Hortaliza"
